$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns into the existing data block:
#   - one before the old "metalAlloys" column (F) -> becomes the "water" column
#   - one before the old "metalAlloysWaste" column (which, after the first
#     insert, sits at J) -> becomes the "food" column
$ws.Columns("F").Insert()
$ws.Columns("J").Insert()

# Header row labels for the two inserted columns, plus a brand new trailing
# column "foodWaste" appended after the existing data (N).
$ws.Range("F1").Value = "water"
$ws.Range("J1").Value = "food"
$ws.Range("N1").Value = "foodWaste"
$ws.Range("N1").HorizontalAlignment = -4108

# Fill in the values for the new / now-uniform resource columns (F:N) for
# every data row (2-7) with 100.
for ($r = 2; $r -le 7; $r++) {
    for ($c = 6; $c -le 14; $c++) {
        $ws.Cells.Item($r, $c).Value = 100
    }
}

# Match the final selection/active cell shown in the workbook.
$ws.Range("J10").Select()

# The saved workbook also carries an explicit (portrait) page setup.
$ws.PageSetup.Orientation = 1
